$wb = $excel.ActiveWorkbook

# --- Rename the second sheet: "Include from ActClass" -> "Include #0" ---
$ws2 = $wb.Worksheets.Item("Include from ActClass")
$ws2.Name = "Include #0"

# --- Update the Metadata sheet ---
$ws = $wb.Worksheets.Item("Metadata")

# Version bump and regeneration date
$ws.Range("B3").Value = "2.0.1-sd-202510-matchbox-patch"
$ws.Range("B8").Value = "2025-10-29T22:15:57+01:00"

# Insert a new "Jurisdiction" property row after "Contact" (row 10) and
# before "Description" (row 11), pushing Description/Purpose/Copyright/
# Immutable down by one row. Copy row 14's formatting onto the
# newly-needed row 15 first so every row keeps the same cell style
# as its neighbours once the values are shifted down.
$ws.Range("A14:B14").Copy()
$ws.Range("A15:B15").PasteSpecial(-4122)

$ws.Range("A15").Value = $ws.Range("A14").Value2
$ws.Range("B15").Value = $ws.Range("B14").Value2

$ws.Range("A14").Value = $ws.Range("A13").Value2
$ws.Range("B14").Value = $ws.Range("B13").Value2

$ws.Range("A13").Value = $ws.Range("A12").Value2
$ws.Range("B13").Value = $ws.Range("B12").Value2

$ws.Range("A12").Value = $ws.Range("A11").Value2
$ws.Range("B12").Value = $ws.Range("B11").Value2

$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

Write-Output "done"
